# "Dumping code skeleton done"
#
# The "artist_*" field names in column A are renamed to drop the
# "artist_" prefix (artist_id -> id, artist_playmeid -> playmeid, ...).
# A few of the new, shorter names (latitude / longitude) already exist
# elsewhere in the shared-strings table, so Excel will naturally reuse
# those entries once the old "artist_*" strings become orphaned and are
# dropped from xl/sharedStrings.xml on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Artist identity" block (rows 2-5)
$ws.Range("A2").Value = "id"
$ws.Range("A3").Value = "playmeid"
$ws.Range("A4").Value = "digi_7id"
$ws.Range("A5").Value = "mbid"

# "Artist info" block (rows 7-14)
$ws.Range("A7").Value = "id"
$ws.Range("A8").Value = "mbid"
$ws.Range("A9").Value = "name"
$ws.Range("A10").Value = "familiarity"
$ws.Range("A11").Value = "hottness"
$ws.Range("A12").Value = "latitude"
$ws.Range("A13").Value = "location"
$ws.Range("A14").Value = "longitude"

# Foreign-key summary block near the bottom (rows 47-48)
$ws.Range("A47").Value = "id"
$ws.Range("A48").Value = "mbid"

# Row 44 ("time_signature_confidence") ends up taller after the edits
# (its row grows to fit the wrapped label).
$ws.Rows.Item(44).RowHeight = 28.2

# Active selection moves from C8 to C5.
$ws.Range("C5").Select()
